$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sending-cluster label "Resolving-Mac" -> "Neutrophils" (rows 10-13, col A).
# (Target-cluster column D already reads "Neutrophils" for those rows, so no visible
# change is needed there even though the underlying shared-string slot is consolidated.)
$ws.Range("A10:A13").Value = "Neutrophils"

# Recomputed NATMI TPM-derived metrics (ligand/receptor expression, specificity and
# edge-weight columns) for every sending/target-cluster pair.
$cellUpdates = @(
    @{ Cell = "G2"; Value = 0.4815393333333333 },
    @{ Cell = "H2"; Value = 1.444618 },
    @{ Cell = "I2"; Value = 0.3617347224948818 },
    @{ Cell = "J2"; Value = 0.3617347224948818 },
    @{ Cell = "M2"; Value = 1.139366 },
    @{ Cell = "N2"; Value = 3.418098 },
    @{ Cell = "O2"; Value = 0.1546713947032042 },
    @{ Cell = "P2"; Value = 0.1546713947032042 },
    @{ Cell = "Q2"; Value = 0.5486495440626667 },
    @{ Cell = "R2"; Value = 4.937845896564 },
    @{ Cell = "S2"; Value = 0.05595001404085988 },
    @{ Cell = "T2"; Value = 0.05595001404085988 },
    @{ Cell = "G3"; Value = 0.4815393333333333 },
    @{ Cell = "H3"; Value = 1.444618 },
    @{ Cell = "I3"; Value = 0.3617347224948818 },
    @{ Cell = "J3"; Value = 0.3617347224948818 },
    @{ Cell = "O3"; Value = 0.7676983257595695 },
    @{ Cell = "P3"; Value = 0.7676983257595695 },
    @{ Cell = "Q3"; Value = 2.723175395255778 },
    @{ Cell = "R3"; Value = 24.508578557302 },
    @{ Cell = "S3"; Value = 0.2777031408284232 },
    @{ Cell = "T3"; Value = 0.2777031408284232 },
    @{ Cell = "G4"; Value = 0.4815393333333333 },
    @{ Cell = "H4"; Value = 1.444618 },
    @{ Cell = "I4"; Value = 0.3617347224948818 },
    @{ Cell = "J4"; Value = 0.3617347224948818 },
    @{ Cell = "M4"; Value = 0.4872916666666667 },
    @{ Cell = "N4"; Value = 1.461875 },
    @{ Cell = "O4"; Value = 0.06615089594615092 },
    @{ Cell = "P4"; Value = 0.06615089594615092 },
    @{ Cell = "Q4"; Value = 0.2346501043055556 },
    @{ Cell = "R4"; Value = 2.11185093875 },
    @{ Cell = "S4"; Value = 0.0239290759878687 },
    @{ Cell = "T4"; Value = 0.0239290759878687 },
    @{ Cell = "G5"; Value = 0.4815393333333333 },
    @{ Cell = "H5"; Value = 1.444618 },
    @{ Cell = "I5"; Value = 0.3617347224948818 },
    @{ Cell = "J5"; Value = 0.3617347224948818 },
    @{ Cell = "M5"; Value = 0.08456133333333334 },
    @{ Cell = "N5"; Value = 0.253684 },
    @{ Cell = "O5"; Value = 0.0114793835910754 },
    @{ Cell = "P5"; Value = 0.0114793835910754 },
    @{ Cell = "Q5"; Value = 0.04071960807911111 },
    @{ Cell = "R5"; Value = 0.366476472712 },
    @{ Cell = "S5"; Value = 0.004152491637729959 },
    @{ Cell = "T5"; Value = 0.00415249163772996 },
    @{ Cell = "E6"; Value = 1 },
    @{ Cell = "F6"; Value = 0.3333333333333333 },
    @{ Cell = "G6"; Value = 0.259826 },
    @{ Cell = "H6"; Value = 0.779478 },
    @{ Cell = "I6"; Value = 0.1951825728468463 },
    @{ Cell = "J6"; Value = 0.1951825728468463 },
    @{ Cell = "M6"; Value = 1.139366 },
    @{ Cell = "N6"; Value = 3.418098 },
    @{ Cell = "O6"; Value = 0.1546713947032042 },
    @{ Cell = "P6"; Value = 0.1546713947032042 },
    @{ Cell = "Q6"; Value = 0.296036910316 },
    @{ Cell = "R6"; Value = 2.664332192844 },
    @{ Cell = "S6"; Value = 0.03018916076398147 },
    @{ Cell = "T6"; Value = 0.03018916076398147 },
    @{ Cell = "E7"; Value = 1 },
    @{ Cell = "F7"; Value = 0.3333333333333333 },
    @{ Cell = "G7"; Value = 0.259826 },
    @{ Cell = "H7"; Value = 0.779478 },
    @{ Cell = "I7"; Value = 0.1951825728468463 },
    @{ Cell = "J7"; Value = 0.1951825728468463 },
    @{ Cell = "O7"; Value = 0.7676983257595695 },
    @{ Cell = "P7"; Value = 0.7676983257595695 },
    @{ Cell = "Q7"; Value = 1.469354051204667 },
    @{ Cell = "R7"; Value = 13.224186460842 },
    @{ Cell = "S7"; Value = 0.1498413343919691 },
    @{ Cell = "T7"; Value = 0.1498413343919691 },
    @{ Cell = "E8"; Value = 1 },
    @{ Cell = "F8"; Value = 0.3333333333333333 },
    @{ Cell = "G8"; Value = 0.259826 },
    @{ Cell = "H8"; Value = 0.779478 },
    @{ Cell = "I8"; Value = 0.1951825728468463 },
    @{ Cell = "J8"; Value = 0.1951825728468463 },
    @{ Cell = "M8"; Value = 0.4872916666666667 },
    @{ Cell = "N8"; Value = 1.461875 },
    @{ Cell = "O8"; Value = 0.06615089594615092 },
    @{ Cell = "P8"; Value = 0.06615089594615092 },
    @{ Cell = "Q8"; Value = 0.1266110445833333 },
    @{ Cell = "R8"; Value = 1.13949940125 },
    @{ Cell = "S8"; Value = 0.01291150206689375 },
    @{ Cell = "T8"; Value = 0.01291150206689375 },
    @{ Cell = "E9"; Value = 1 },
    @{ Cell = "F9"; Value = 0.3333333333333333 },
    @{ Cell = "G9"; Value = 0.259826 },
    @{ Cell = "H9"; Value = 0.779478 },
    @{ Cell = "I9"; Value = 0.1951825728468463 },
    @{ Cell = "J9"; Value = 0.1951825728468463 },
    @{ Cell = "M9"; Value = 0.08456133333333334 },
    @{ Cell = "N9"; Value = 0.253684 },
    @{ Cell = "O9"; Value = 0.0114793835910754 },
    @{ Cell = "P9"; Value = 0.0114793835910754 },
    @{ Cell = "Q9"; Value = 0.02197123299466667 },
    @{ Cell = "R9"; Value = 0.197741096952 },
    @{ Cell = "S9"; Value = 0.002240575624001967 },
    @{ Cell = "T9"; Value = 0.002240575624001967 },
    @{ Cell = "G10"; Value = 0.5898293333333333 },
    @{ Cell = "H10"; Value = 1.769488 },
    @{ Cell = "I10"; Value = 0.4430827046582718 },
    @{ Cell = "J10"; Value = 0.4430827046582719 },
    @{ Cell = "M10"; Value = 1.139366 },
    @{ Cell = "N10"; Value = 3.418098 },
    @{ Cell = "O10"; Value = 0.1546713947032042 },
    @{ Cell = "P10"; Value = 0.1546713947032042 },
    @{ Cell = "Q10"; Value = 0.6720314882026667 },
    @{ Cell = "R10"; Value = 6.048283393824 },
    @{ Cell = "S10"; Value = 0.06853221989836279 },
    @{ Cell = "T10"; Value = 0.06853221989836281 },
    @{ Cell = "G11"; Value = 0.5898293333333333 },
    @{ Cell = "H11"; Value = 1.769488 },
    @{ Cell = "I11"; Value = 0.4430827046582718 },
    @{ Cell = "J11"; Value = 0.4430827046582719 },
    @{ Cell = "O11"; Value = 0.7676983257595695 },
    @{ Cell = "P11"; Value = 0.7676983257595695 },
    @{ Cell = "Q11"; Value = 3.335571191692444 },
    @{ Cell = "R11"; Value = 30.020140725232 },
    @{ Cell = "S11"; Value = 0.3401538505391771 },
    @{ Cell = "T11"; Value = 0.3401538505391771 },
    @{ Cell = "G12"; Value = 0.5898293333333333 },
    @{ Cell = "H12"; Value = 1.769488 },
    @{ Cell = "I12"; Value = 0.4430827046582718 },
    @{ Cell = "J12"; Value = 0.4430827046582719 },
    @{ Cell = "M12"; Value = 0.4872916666666667 },
    @{ Cell = "N12"; Value = 1.461875 },
    @{ Cell = "O12"; Value = 0.06615089594615092 },
    @{ Cell = "P12"; Value = 0.06615089594615092 },
    @{ Cell = "Q12"; Value = 0.2874189188888889 },
    @{ Cell = "R12"; Value = 2.58677027 },
    @{ Cell = "S12"; Value = 0.02931031789138846 },
    @{ Cell = "T12"; Value = 0.02931031789138846 },
    @{ Cell = "G13"; Value = 0.5898293333333333 },
    @{ Cell = "H13"; Value = 1.769488 },
    @{ Cell = "I13"; Value = 0.4430827046582718 },
    @{ Cell = "J13"; Value = 0.4430827046582719 },
    @{ Cell = "M13"; Value = 0.08456133333333334 },
    @{ Cell = "N13"; Value = 0.253684 },
    @{ Cell = "O13"; Value = 0.0114793835910754 },
    @{ Cell = "P13"; Value = 0.0114793835910754 },
    @{ Cell = "Q13"; Value = 0.04987675486577778 },
    @{ Cell = "R13"; Value = 0.448890793792 },
    @{ Cell = "S13"; Value = 0.005086316329343473 },
    @{ Cell = "T13"; Value = 0.005086316329343475 }
)

foreach ($u in $cellUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
